$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1) - copy the existing header formatting
# (bold, bordered, centered) from the adjacent "sum" header, then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New "Save" column data rows (H2:H3), unstyled like the other data cells.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
